# "Generate Report for Handback"
#
# The localization-status report is refreshed with handback information for
# both locale sheets (zh-cn, de-de):
#   - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#   - The "Latest Target File" (F) and "Latest Handback File" (G) columns are
#     populated with links to the handed-back files
#   - The "Latest Handback DateTime" (H) column is stamped with the handback time

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2 - 203746f0-9d87-4f07-9364-eebc27e565a5
$wsZh.Range("C2").Value = $statusText
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/7a473f2b05082f979479e2f5043ea50fd3b4076d/e2e/203746f0-9d87-4f07-9364-eebc27e565a5.md",
    [Type]::Missing,
    [Type]::Missing,
    "203746f0-9d87-4f07-9364-eebc27e565a5.md") | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/39cd67740c3ef0b5bdfa0a2cf0c4894f34cb96fa/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/203746f0-9d87-4f07-9364-eebc27e565a5.b1d64b89cfbac632499873b3d217eb7fc4c4a24d.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "203746f0-9d87-4f07-9364-eebc27e565a5.b1d64b89cfbac632499873b3d217eb7fc4c4a24d.zh-cn.xlf") | Out-Null
$wsZh.Range("H2").Value = "2016-03-13 02:57:15"

# Row 3 - b4b85e53-1c06-483f-81d4-36d4c4d2cf74
$wsZh.Range("C3").Value = $statusText
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/7a473f2b05082f979479e2f5043ea50fd3b4076d/e2e/b4b85e53-1c06-483f-81d4-36d4c4d2cf74.md",
    [Type]::Missing,
    [Type]::Missing,
    "b4b85e53-1c06-483f-81d4-36d4c4d2cf74.md") | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/39cd67740c3ef0b5bdfa0a2cf0c4894f34cb96fa/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/b4b85e53-1c06-483f-81d4-36d4c4d2cf74.ed9d9b2555ca648603c3c407164ab36cd5198600.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "b4b85e53-1c06-483f-81d4-36d4c4d2cf74.ed9d9b2555ca648603c3c407164ab36cd5198600.zh-cn.xlf") | Out-Null
$wsZh.Range("H3").Value = "2016-03-13 02:57:15"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 2 - 203746f0-9d87-4f07-9364-eebc27e565a5
$wsDe.Range("C2").Value = $statusText
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/7a473f2b05082f979479e2f5043ea50fd3b4076d/e2e/203746f0-9d87-4f07-9364-eebc27e565a5.md",
    [Type]::Missing,
    [Type]::Missing,
    "203746f0-9d87-4f07-9364-eebc27e565a5.md") | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/94ff035f2bd545cc61a5993ba13fa65d05639864/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/203746f0-9d87-4f07-9364-eebc27e565a5.b1d64b89cfbac632499873b3d217eb7fc4c4a24d.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "203746f0-9d87-4f07-9364-eebc27e565a5.b1d64b89cfbac632499873b3d217eb7fc4c4a24d.de-de.xlf") | Out-Null
$wsDe.Range("H2").Value = "2016-03-13 02:57:22"

# Row 3 - b4b85e53-1c06-483f-81d4-36d4c4d2cf74
$wsDe.Range("C3").Value = $statusText
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/7a473f2b05082f979479e2f5043ea50fd3b4076d/e2e/b4b85e53-1c06-483f-81d4-36d4c4d2cf74.md",
    [Type]::Missing,
    [Type]::Missing,
    "b4b85e53-1c06-483f-81d4-36d4c4d2cf74.md") | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/94ff035f2bd545cc61a5993ba13fa65d05639864/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/b4b85e53-1c06-483f-81d4-36d4c4d2cf74.ed9d9b2555ca648603c3c407164ab36cd5198600.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "b4b85e53-1c06-483f-81d4-36d4c4d2cf74.ed9d9b2555ca648603c3c407164ab36cd5198600.de-de.xlf") | Out-Null
$wsDe.Range("H3").Value = "2016-03-13 02:57:22"

Write-Host "Localization status report updated for handback."
